$p = $ppt.ActivePresentation
$s = $p.Slides.Item(24)
$shp = $s.Shapes.Item(2)
$tf = $shp.TextFrame
$tr = $tf.TextRange

# 1. Merge the split "https://javadoc.io/doc/org.mockito/" + "m" +
#    "ockito-core/latest/org/mockito/Mockito.html" runs back into a single run.
$mockitoPara = $tr.Paragraphs(2)
$mockitoPara.Text = "TEMP_PLACEHOLDER"
$mockitoPara.Text = "https://javadoc.io/doc/org.mockito/mockito-core/latest/org/mockito/Mockito.html"

# 2. Add a new paragraph before the junit.org reference line, linking to the
#    GitHub repo (reusing the junit.org hyperlink's run formatting).
$junitPara = $tr.Paragraphs(1)
[void]$junitPara.InsertBefore("https://github.com/gitaroktato/junit5-lectures`r")
